$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so decimal-looking strings (e.g. "1.004") are not
# auto-coerced to numbers; we restore the Normal style afterwards so cells end up
# without an explicit style override, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.019.76"
$ws.Range("E2").Value = "  -1.78%  "
$ws.Range("D3").Value = "1.790.46"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "313.94"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "0.5219"
$ws.Range("E7").Value = "  +2.24%  "
$ws.Range("D8").Value = "0.3802"
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("D9").Value = "0.07915"
$ws.Range("E9").Value = "  -3.15%  "
$ws.Range("D10").Value = "41.35"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "1.090"
$ws.Range("E11").Value = "  -1.87%  "
$ws.Range("D12").Value = "6.265"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "1.004"
$ws.Range("E13").Value = "  +0.29%  "
$ws.Range("D14").Value = "20.44"
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("D15").Value = "1.794.59"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("D16").Value = "7.271"
$ws.Range("E16").Value = "  -3.41%  "
$ws.Range("D17").Value = "91.95"
$ws.Range("D18").Value = "0.00001082"
$ws.Range("E18").Value = "  -4.09%  "
$ws.Range("D19").Value = "0.06531"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "17.24"
$ws.Range("E21").Value = "  -3.26%  "
$ws.Range("D22").Value = "5.938"
$ws.Range("E22").Value = "  -2.58%  "
$ws.Range("D23").Value = "28.063.24"
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("D24").Value = "11.10"
$ws.Range("E24").Value = "  -2.77%  "
$ws.Range("D25").Value = "2.260"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "160.98"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").Value = "20.38"
$ws.Range("E27").Value = "  -5.14%  "
$ws.Range("D28").Value = "1.994.68"
$ws.Range("E28").Value = "  -1.93%  "
$ws.Range("D29").Value = "2.316"
$ws.Range("E29").Value = "  -3.43%  "
$ws.Range("D30").Value = "122.32"
$ws.Range("E30").Value = "  -3.35%  "
$ws.Range("D31").Value = "0.1073"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "1.046"
$ws.Range("E32").Value = "  -5.74%  "
$ws.Range("D33").Value = "3.674"
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("D34").Value = "5.518"
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("D35").Value = "0.07241"
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("D36").Value = "12.18"
$ws.Range("E36").Value = "  +7.74%  "
$ws.Range("D37").Value = "0.02315"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").Value = "8.759"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").Value = "0.2130"
$ws.Range("E39").Value = "  -4.46%  "
$ws.Range("D40").Value = "5.045"
$ws.Range("E40").Value = "  -4.34%  "
$ws.Range("D41").Value = "0.6130"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").Value = "1.164"
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").Value = "13.26"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "3.764"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "0.5921"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "127.84"
$ws.Range("E47").Value = "  +1.92%  "
$ws.Range("D48").Value = "1.226"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("D50").Value = "0.06728"
$ws.Range("E50").Value = "  -2.77%  "
$ws.Range("D51").Value = "72.67"
$ws.Range("E51").Value = "  -1.80%  "

# Restore default styling on column D (clears the temporary text format).
$ws.Range("D2:D51").Style = "Normal"
